{"js": "// Add a new \"To do\" bullet right after the existing\n// \"Checken of alle routes naar tabnaam/pagenaam gaan OK\" item:\n// a bold BodyText paragraph reading\n// \"custom nav nog toevoegen aan gameedit en userstats\".\n\nconst body = context.document.body;\n\n// Locate the paragraph that ends the \"Checken of alle routes...\" entry so\n// the new paragraph is inserted directly after it (and before the\n// \"GameAddEdit: datum verwerken\" paragraph that currently follows it).\nconst searchResults = body.search(\n  \"Checken of alle routes naar tabnaam/pagenaam gaan OK\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\n    \"Could not find the 'Checken of alle routes...' paragraph to anchor the new note after.\"\n  );\n}\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// Insert the new paragraph after the anchor; it inherits the BodyText\n// style/spacing from the paragraph it is split from.\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"custom nav nog toevoegen aan gameedit en userstats\",\n  Word.InsertLocation.after\n);\n\n// Make sure the whole new line is bold (matches the other bold \"labels\" /\n// \"OK\" markers used throughout this to-do list).\nnewParagraph.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Add a new \"To do\" bullet right after the existing\n# \"Checken of alle routes naar tabnaam/pagenaam gaan OK\" item:\n# a bold BodyText paragraph reading\n# \"custom nav nog toevoegen aan gameedit en userstats\".\n\n$d = $word.ActiveDocument\n$anchorText = \"Checken of alle routes naar tabnaam/pagenaam gaan OK\"\n\n# Locate the paragraph to anchor the new note after (Range.Text carries a\n# trailing paragraph mark, so trim before comparing).\n$anchorIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq 0) {\n    throw \"Could not find the 'Checken of alle routes...' paragraph to anchor the new note after.\"\n}\n\n# Insert an empty paragraph right after the anchor paragraph; it inherits\n# the BodyText style/spacing from the paragraph it is split from.\n$anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n$anchorRange.Collapse(0)  # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n\n# Fill in the new paragraph's text and make sure the whole line is bold\n# (matches the other bold \"labels\" / \"OK\" markers used in this to-do list).\n$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n$newParagraph.Range.Text = \"custom nav nog toevoegen aan gameedit en userstats\"\n$newParagraph.Range.Font.Bold = 1\n$newParagraph.Range.Font.BoldBi = 1\n"}
